$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 142 ---
$ws.Cells.Item(141, 1).Copy()
$ws.Cells.Item(142, 1).PasteSpecial(-4122)
$ws.Cells.Item(142, 1).Value = 45491.2916666667

$ws.Cells.Item(142, 2).Value = 0
$ws.Cells.Item(142, 3).Value = 2.25
$ws.Cells.Item(142, 4).Value = 2.25
$ws.Cells.Item(142, 5).Value = 2.25
$ws.Cells.Item(142, 6).Value = 2.25

$ws.Cells.Item(142, 7).NumberFormat = "@"
$ws.Cells.Item(142, 7).Value = "2.25"
$ws.Cells.Item(142, 7).Style = "Normal"

$ws.Cells.Item(142, 8).Value = "LS.MI"

# --- Row 143 ---
$ws.Cells.Item(141, 1).Copy()
$ws.Cells.Item(143, 1).PasteSpecial(-4122)
$ws.Cells.Item(143, 1).Value = 45492.3508449074

$ws.Cells.Item(143, 2).Value = 9000
$ws.Cells.Item(143, 3).Value = 2.33999991416931
$ws.Cells.Item(143, 4).Value = 2.25999999046326
$ws.Cells.Item(143, 5).Value = 2.28999996185303
$ws.Cells.Item(143, 6).Value = 2.25999999046326

$ws.Cells.Item(143, 7).NumberFormat = "@"
$ws.Cells.Item(143, 7).Value = "2.25999999046326"
$ws.Cells.Item(143, 7).Style = "Normal"

$ws.Cells.Item(143, 8).Value = "LS.MI"
